# Updates cryptocurrency price/volume figures in the "cryptos" worksheet
# to reflect the latest scrape (GitHub Actions symbol-list refresh).
# Values in columns D (Price) and E (Volume/1h) are stored as plain text
# (e.g. "330.39", "0.08%"), so a leading apostrophe is used to force
# Excel to keep them as text instead of auto-converting to numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "'330.39"
$ws.Range("E2").Value = "'0.08%"
$ws.Range("D3").Value = "'41.63"
$ws.Range("E3").Value = "'1.12%"
$ws.Range("D4").Value = "'5.703"
$ws.Range("E4").Value = "'0.13%"
$ws.Range("D5").Value = "'0.08425"
$ws.Range("E5").Value = "'4.41%"
$ws.Range("D6").Value = "'8.826"
$ws.Range("E6").Value = "'1.15%"
$ws.Range("D7").Value = "'1.991"
$ws.Range("E7").Value = "'-1.92%"
$ws.Range("D8").Value = "'4.488"
$ws.Range("E8").Value = "'-0.59%"
$ws.Range("D9").Value = "'2.923"
$ws.Range("E9").Value = "'-0.40%"
$ws.Range("D10").Value = "'0.9266"
$ws.Range("E10").Value = "'0.44%"
$ws.Range("D11").Value = "'0.1280"
$ws.Range("E11").Value = "'0.73%"
$ws.Range("D12").Value = "'0.1963"
$ws.Range("E12").Value = "'1.11%"
$ws.Range("D13").Value = "'0.09306"
$ws.Range("E13").Value = "'0.09%"
$ws.Range("D14").Value = "'0.03951"
$ws.Range("E14").Value = "'6.63%"
$ws.Range("E15").Value = "'0.91%"
$ws.Range("D16").Value = "'0.001305"
$ws.Range("E16").Value = "'-0.01%"
$ws.Range("D17").Value = "'0.006116"
$ws.Range("E17").Value = "'-2.22%"
$ws.Range("D18").Value = "'3.431"
$ws.Range("E18").Value = "'1.96%"
$ws.Range("E19").Value = "'0.72%"
$ws.Range("D20").Value = "'8.956"
$ws.Range("E20").Value = "'8.63%"
$ws.Range("D21").Value = "'0.1365"
$ws.Range("E21").Value = "'-3.68%"
$ws.Range("D22").Value = "'0.2513"
$ws.Range("E22").Value = "'-5.30%"
$ws.Range("D23").Value = "'0.04423"
$ws.Range("E23").Value = "'-0.15%"
$ws.Range("D24").Value = "'0.001247"
$ws.Range("E24").Value = "'-1.02%"
$ws.Range("D25").Value = "'0.004374"
$ws.Range("E25").Value = "'0.83%"
$ws.Range("D26").Value = "'0.0001192"
$ws.Range("E26").Value = "'-4.01%"
$ws.Range("D27").Value = "'0.0003997"
$ws.Range("E27").Value = "'0.10%"
$ws.Range("D39").Value = "'0.02823"
$ws.Range("E39").Value = "'-1.45%"
$ws.Range("D40").Value = "'0.05519"
$ws.Range("E40").Value = "'0.94%"
$ws.Range("D41").Value = "'0.007920"
$ws.Range("E41").Value = "'3.96%"
$ws.Range("E42").Value = "'1.54%"
$ws.Range("D43").Value = "'0.008981"
$ws.Range("E43").Value = "'-9.84%"
$ws.Range("D44").Value = "'0.002093"
$ws.Range("E44").Value = "'-1.85%"
$ws.Range("E45").Value = "'-7.07%"
$ws.Range("D46").Value = "'0.00007344"
$ws.Range("E46").Value = "'8.57%"
$ws.Range("E47").Value = "'0.13%"
$ws.Range("D48").Value = "'0.003261"
$ws.Range("E48").Value = "'9.04%"
$ws.Range("D49").Value = "'0.002282"
$ws.Range("E49").Value = "'0.06%"
$ws.Range("D50").Value = "'0.00002104"
$ws.Range("E50").Value = "'0.13%"
$ws.Range("E51").Value = "'0.13%"
